# Commit: "this is new modify code for excel update in 14/06/2024"
#
# 1) Fix a couple of mobile numbers that were mistyped.
# 2) Flip the "Result output" message on row 2 back to the normal
#    "OrangeContact number is already in use. apple" message.
# 3) Re-key the owner on rows 3/4 ("PRASANNA TARAI") down to shorter
#    name variants and log three more validation attempts as new rows
#    (5, 6, 7) underneath the existing data.
# 4) Leave the selection on H2 (where the last edit happened).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix mistyped mobile numbers -------------------------------------------------
$ws.Range("H2").Value = 6372602823
$ws.Range("H4").Value = 8328961074

# --- row 2's logged result goes back to the generic "apple" message -------------
$ws.Range("AF2").Value = "OrangeContact number is already in use. apple"

# --- append new rows 5-7, copying formatting from row 4 so styles match ---------
$ws.Range("A4:AF4").Copy()
$ws.Range("A5:AF5").PasteSpecial(-4122)
$ws.Range("A4:AF4").Copy()
$ws.Range("A6:AF6").PasteSpecial(-4122)
$ws.Range("A4:AF4").Copy()
$ws.Range("A7:AF7").PasteSpecial(-4122)

# Row 5
$ws.Range("A5").Value = 6405
$ws.Range("B5").Value = "PRASANN"
$ws.Range("C5").Value = "TARA"
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 112
$ws.Range("F5").Value = 234
$ws.Range("G5").Value = "SINGITALIA"
$ws.Range("H5").Value = 8324961274
$ws.Range("I5").Value = "rajatmohantysahapur@gmail.com"
$ws.Range("AF5").Value = "OrangeContact number is already in use. apple"

# Row 6
$ws.Range("A6").Value = 6405
$ws.Range("B6").Value = "PRASANN"
$ws.Range("C6").Value = "TARA"
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 112
$ws.Range("F6").Value = 234
$ws.Range("G6").Value = "SINGITALIA"
$ws.Range("H6").Value = 8324961274
$ws.Range("I6").Value = "rajatmohantysahapur@gmail.com"
$ws.Range("AF6").Value = "OrangeContact number is already in use. apple"

# Row 7
$ws.Range("A7").Value = 6405
$ws.Range("B7").Value = "PRASA"
$ws.Range("C7").Value = "TAR"
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = 112
$ws.Range("F7").Value = 234
$ws.Range("G7").Value = "SINGITALIA"
$ws.Range("H7").Value = 8324960276
$ws.Range("I7").Value = "rajatmohantysahapur@gmail.com"
$ws.Range("AF7").Value = "OrangeContact number is already in use. apple"

# --- leave the selection on H2 (matches the author's final cursor position) -----
[void]$ws.Range("H2").Select()
